$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the hazard-level labels in column D (the "requirement issue" fix):
#   "轻度火灾危险" -> "轻度"   (mild)
#   "中度火灾危险" -> "中度"   (moderate)
#   "高度火灾危险" -> "高度"   (high)
#   "严重火灾危险" -> "严重"   (severe)
$ws.Range("D2").Value = "轻度"
$ws.Range("D3").Value = "中度"
$ws.Range("D4").Value = "高度"
$ws.Range("D5").Value = "严重"

# Move / leave the active selection on D5, matching the resaved workbook state.
$ws.Range("D5").Select()
